$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$simRows = @(2, 5, 6, 12, 20, 24, 28, 29, 30, 31)

$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
foreach ($r in $simRows) {
    $ws.Range("H$r").PasteSpecial(-4122)  # xlPasteFormats
}

$ws.Range("H1").Value = "exibir_ao_iniciar"
foreach ($r in $simRows) {
    $ws.Cells.Item($r, 8).Value = "sim"
}
